$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 10 blank rows before row 251 so the existing block (rows 251-300)
# shifts down to rows 261-310.
$ws.Range("A251:C260").EntireRow.Insert()

# Populate the newly inserted key/value pair at row 247 (continuing the
# y2001..y2006 sequence already present in rows 241-246), matching the
# formatting of the row right above it.
$ws.Range("B246").Copy()
$ws.Range("B247").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A247").Value = "y2007"
$ws.Range("B247").Value = '주어진 조건에서 함수의 개형을 파악내고 극대와 극소의 $x$ 좌표를 알아내서 도함수의 식을 통해 $f(x)$를 표현합니다.'

# Match the author's final viewport / selection position recorded in the diff.
$ws.Application.ActiveWindow.ScrollRow = 232
$ws.Range("B250").Select()
